$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A (Opportunity ID) to text so numeric-looking IDs stay as strings
$ws.Range("A2:A16").NumberFormat = "@"

$ws.Range("A2").Value = "1327168"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1327168"
$ws.Range("C2").Value = "AI & ML Intern"
$ws.Range("D2").Value = "Manipal, Karnataka, India"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "2 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "M.A.H.E."

$ws.Range("A3").Value = "1326776"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1326776"
$ws.Range("C3").Value = "Biotechnology Intern"
$ws.Range("D3").Value = "Manipal, Karnataka, India"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("G3").Value = "3 - 6 Months"
$ws.Range("H3").Value = "M.A.H.E."

$ws.Range("A4").Value = "1326767"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1326767"
$ws.Range("C4").Value = "Machine Learning Intern"
$ws.Range("D4").Value = "Manipal, Karnataka, India"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "M.A.H.E."

$ws.Range("A5").Value = "1326765"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1326765"
$ws.Range("C5").Value = "Electrical Engineering Intern"
$ws.Range("D5").Value = "Manipal, Karnataka, India"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "0 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "M.A.H.E."

$ws.Range("A6").Value = "1326761"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1326761"
$ws.Range("C6").Value = "Civil Engineering Intern"
$ws.Range("D6").Value = "Manipal, Karnataka, India"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "0 applicants"
$ws.Range("G6").Value = "9 - 12 Weeks"
$ws.Range("H6").Value = "M.A.H.E."

$ws.Range("A7").Value = "1326760"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1326760"
$ws.Range("C7").Value = "Architecture Intern"
$ws.Range("D7").Value = "Manipal, Karnataka, India"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "0 applicants"
$ws.Range("G7").Value = "9 - 12 Weeks"
$ws.Range("H7").Value = "M.A.H.E."

$ws.Range("A8").Value = "1326757"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1326757"
$ws.Range("C8").Value = "Architectural Intern"
$ws.Range("D8").Value = "Manipal, Karnataka, India"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "0 applicants"
$ws.Range("G8").Value = "9 - 12 Weeks"
$ws.Range("H8").Value = "M.A.H.E."

$ws.Range("A9").Value = "1326756"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1326756"
$ws.Range("C9").Value = "Electrical & ML Intern"
$ws.Range("D9").Value = "Manipal, Karnataka, India"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "0 applicants"
$ws.Range("G9").Value = "9 - 12 Weeks"
$ws.Range("H9").Value = "M.A.H.E."

$ws.Range("A10").Value = "1326666"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1326666"
$ws.Range("C10").Value = "Computer Engineering Intern"
$ws.Range("D10").Value = "Manipal, Karnataka, India"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "0 applicants"
$ws.Range("G10").Value = "9 - 12 Weeks"
$ws.Range("H10").Value = "M.A.H.E."

$ws.Range("A11").Value = "1326661"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1326661"
$ws.Range("C11").Value = "AL & ML Intern"
$ws.Range("D11").Value = "Manipal, Karnataka, India"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "0 applicants"
$ws.Range("G11").Value = "9 - 12 Weeks"
$ws.Range("H11").Value = "M.A.H.E."

$ws.Range("A12").Value = "1326639"
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1326639"
$ws.Range("C12").Value = "Machine Learning Intern"
$ws.Range("D12").Value = "Manipal, Karnataka, India"
$ws.Range("E12").Value = "No"
$ws.Range("F12").Value = "0 applicants"
$ws.Range("G12").Value = "9 - 12 Weeks"
$ws.Range("H12").Value = "M.A.H.E."

$ws.Range("A13").Value = "1326094"
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1326094"
$ws.Range("C13").Value = "Public Policy, Government Affairs Support Intern"
$ws.Range("D13").Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Range("E13").Value = "No"
$ws.Range("F13").Value = "11 applicants"
$ws.Range("G13").Value = "6 - 18 Months"
$ws.Range("H13").Value = "Samsung Electronics Latinoamerica SELA"

$ws.Range("A14").Value = "1317568"
$ws.Range("B14").Value = "https://aiesec.org/opportunity/global-talent/1317568"
$ws.Range("C14").Value = "Research Intern, Project base learning"
$ws.Range("D14").Value = "Pune, Maharashtra, India"
$ws.Range("E14").Value = "No"
$ws.Range("F14").Value = "20 applicants"
$ws.Range("G14").Value = "9 - 12 Weeks"
$ws.Range("H14").Value = "PVG’s College of Engineering and Technology"

$ws.Range("A15").Value = "1317539"
$ws.Range("B15").Value = "https://aiesec.org/opportunity/global-talent/1317539"
$ws.Range("C15").Value = "HR Intern"
$ws.Range("D15").Value = "Pune, Maharashtra, India"
$ws.Range("E15").Value = "No"
$ws.Range("F15").Value = "9 applicants"
$ws.Range("G15").Value = "6 - 18 Months"
$ws.Range("H15").Value = "Forbes Marshall Pvt Ltd"

$ws.Range("A16").Value = "1317364"
$ws.Range("B16").Value = "https://aiesec.org/opportunity/global-talent/1317364"
$ws.Range("C16").Value = "Business Development cum Marketing Executive"
$ws.Range("D16").Value = "Simpang Ampat, Penang, Malaisie"
$ws.Range("E16").Value = "No"
$ws.Range("F16").Value = "130 applicants"
$ws.Range("G16").Value = "6 - 18 Months"
$ws.Range("H16").Value = "Cavalier Capital Holdings Sdn Bhd"

# Update column widths to match new content
$ws.Columns.Item(4).ColumnWidth = 37.103333333333334
$ws.Columns.Item(7).ColumnWidth = 15.103333333333333
$ws.Columns.Item(8).ColumnWidth = 45.10333333333334

